$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up row 6: vendor was "Book" (a data-entry mistake); correct it to
# "Walmart" and fill in the missing Item1/Qty1 columns for Gabriel Alfaro.
$ws.Range("A6").Value = "Walmart"
$ws.Range("D6").Value = "Towel"
$ws.Range("E6").Value = 4

# Move the active selection to G7 (next empty row area for new expense entry).
$null = $ws.Range("G7").Select()
